# Updated cryptos list (price + 1h volume columns) per GitHub Actions scrape.
# Values that "look like" plain numbers (e.g. "1.00", "612.54") must be
# force-stored as TEXT (like the other Price/Volume cells already are),
# otherwise Excel's normal .Value auto-detection would convert them to
# numeric cells and drop significant trailing zeros. We do that by writing
# a quoted-string formula (="1.00"), then Copy + PasteSpecial(values) to
# collapse it back to a static cell without leaving a formula behind and
# without touching cell style (unlike the apostrophe-prefix / NumberFormat
# "@" tricks, which both stamp a new quotePrefix/text style on the cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.420.49'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").Value = '2.660.08'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Formula = "=""1.00"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Formula = "=""612.54"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").Formula = "=""150.34"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +4.31%  '
$ws.Range("D7").Formula = "=""1.00"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Formula = "=""0.390"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +7.63%  '
$ws.Range("D11").Formula = "=""5.63"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Formula = "=""27.92"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").Value = '3.137.87'
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("D15").Value = '64.263.80'
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("E16").Value = '  +2.32%  '
$ws.Range("D17").Value = '2.653.62'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Formula = "=""12.10"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +5.82%  '
$ws.Range("D19").Formula = "=""4.64"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +4.71%  '
$ws.Range("D20").Formula = "=""347.36"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -0.44%  '
$ws.Range("D24").Formula = "=""66.70"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").Formula = "=""1.77"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +15.18%  '
$ws.Range("D26").Formula = "=""1.73"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +5.06%  '
$ws.Range("D27").Formula = "=""9.42"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +8.42%  '
$ws.Range("D28").Formula = "=""562.37"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +2.61%  '
$ws.Range("D29").Formula = "=""8.27"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +4.74%  '
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").Formula = "=""1.00"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  +1.85%  '
$ws.Range("D33").Value = '0.0₃0859'
$ws.Range("E33").Value = '  +6.27%  '
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Formula = "=""5.31"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +3.77%  '
$ws.Range("D36").Formula = "=""168.90"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("D37").Formula = "=""0.409"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").Formula = "=""1.96"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +6.34%  '
$ws.Range("D40").Formula = "=""19.39"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Formula = "=""167.70"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -3.10%  '
$ws.Range("D43").Formula = "=""40.44"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("D44").Formula = "=""3.87"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +2.95%  '
$ws.Range("D45").Formula = "=""0.0577"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Formula = "=""22.08"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Formula = "=""2.01"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +15.13%  '
$ws.Range("D49").Formula = "=""0.0247"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +2.89%  '
$ws.Range("D50").Formula = "=""0.0966"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Formula = "=""19.09"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +1.65%  '
